$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = '''63.756.05'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -1.79%  '
$cell = $ws.Range("D3")
$cell.Value = '''3.137.83'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("E4").Value = '  +0.11%  '
$cell = $ws.Range("D5")
$cell.Value = '''608.86'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.51%  '
$cell = $ws.Range("D6")
$cell.Value = '''146.00'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -4.93%  '
$ws.Range("E7").Value = '  +0.04%  '
$cell = $ws.Range("D8")
$cell.Value = '''3.133.92'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -1.89%  '
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("E10").Value = '  -2.76%  '
$cell = $ws.Range("D11")
$cell.Value = '''5.35'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -4.44%  '
$cell = $ws.Range("D12")
$cell.Value = '''0.470'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -2.07%  '
$ws.Range("E13").Value = '  -2.48%  '
$cell = $ws.Range("D14")
$cell.Value = '''35.37'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -5.03%  '
$cell = $ws.Range("D15")
$cell.Value = '''3.649.25'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.94%  '
$ws.Range("E16").Value = '  +2.53%  '
$cell = $ws.Range("D17")
$cell.Value = '''63.773.79'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.75%  '
$cell = $ws.Range("D18")
$cell.Value = '''3.130.51'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -2.10%  '
$cell = $ws.Range("D19")
$cell.Value = '''6.86'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -3.14%  '
$cell = $ws.Range("D20")
$cell.Value = '''475.17'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -2.14%  '
$cell = $ws.Range("D21")
$cell.Value = '''14.56'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -2.08%  '
$cell = $ws.Range("D22")
$cell.Value = '''0.712'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -1.32%  '
$cell = $ws.Range("D23")
$cell.Value = '''7.97'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +2.06%  '
$cell = $ws.Range("D24")
$cell.Value = '''13.65'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -2.37%  '
$cell = $ws.Range("D25")
$cell.Value = '''83.18'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -2.65%  '
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("E27").Value = '  -4.81%  '
$cell = $ws.Range("D28")
$cell.Value = '''8.46'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -2.76%  '
$cell = $ws.Range("D29")
$cell.Value = '''0.121'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -5.39%  '
$cell = $ws.Range("D30")
$cell.Value = '''7.09'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +1.15%  '
$cell = $ws.Range("D31")
$cell.Value = '''2.09'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -8.91%  '
$ws.Range("E32").Value = '  +0.03%  '
$cell = $ws.Range("D33")
$cell.Value = '''2.69'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -2.05%  '
$cell = $ws.Range("D34")
$cell.Value = '''26.20'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -3.21%  '
$ws.Range("E35").Value = '  +1.18%  '
$cell = $ws.Range("D36")
$cell.Value = '''0.0₃0779'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +6.81%  '
$cell = $ws.Range("D37")
$cell.Value = '''5.98'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.84%  '
$cell = $ws.Range("D38")
$cell.Value = '''52.42'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -4.06%  '
$cell = $ws.Range("D39")
$cell.Value = '''457.67'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -4.22%  '
$cell = $ws.Range("D40")
$cell.Value = '''3.00'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -9.01%  '
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("E42").Value = '  -5.65%  '
$cell = $ws.Range("D43")
$cell.Value = '''8.30'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -3.09%  '
$cell = $ws.Range("D44")
$cell.Value = '''2.859.52'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$cell = $ws.Range("D45")
$cell.Value = '''0.267'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -4.06%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range("D46")
$cell.Value = '''2.29'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -6.10%  '
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D48")
$cell.Value = '''26.35'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -4.70%  '
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Range("D49")
$cell.Value = '''0.999'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  -2.75%  '
$cell = $ws.Range("D51")
$cell.Value = '''118.99'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -1.49%  '

Write-Host "done"